# novos relatorios e atualizacao tabela dados
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Atualizacao da tabela de dados (linhas 7 a 11) ---
# breakwa11
$ws.Range("B7").Value = 308
$ws.Range("C7").Value = 3

# wzxjohn
$ws.Range("B8").Value = 1932
$ws.Range("C8").Value = 37

# chenshaoju
$ws.Range("B9").Value = 42
$ws.Range("C9").Value = 0

# everyx
$ws.Range("B10").Value = 74

# rwasef1830
$ws.Range("B11").Value = 576
$ws.Range("C11").Value = 9

# --- Novos relatorios: marcadores (formatados com sublinhado) ao lado das linhas everyx/rwasef1830 ---
$ws.Range("E10").Font.Underline = 2
$ws.Range("E11").Font.Underline = 2

# --- Configuracao de impressao ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selecao final ---
[void]$ws.Range("E10").Select()
